$d = $word.ActiveDocument

function Append-ConcluidoSuffix($ParaIndex, $OriginalRunsXml) {

    $p = $d.Paragraphs.Item($ParaIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End - 1   # exclude the paragraph mark
    $target = $d.Range($pStart, $pEnd)

    $suffixRuns = @'
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/></w:rPr><w:t>concluíd</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>.</w:t></w:r>
'@

    $xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + $OriginalRunsXml + $suffixRuns + '</w:p><w:sectPr/></w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xmlSnippet)
}

$runs11 = '<w:r w:rsidRPr="003331E2"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Definição da ordem dos Jogadores</w:t></w:r>'

$runs12 = '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Recebimento </w:t></w:r><w:r w:rsidRPr="003331E2"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>e distribuição de Exércitos</w:t></w:r>'

$runs13 = '<w:r w:rsidRPr="003331E2"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ataques</w:t></w:r>'

$runs14 = '<w:r w:rsidRPr="003331E2"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ocupação de território conquistado</w:t></w:r>'

$runs15 = '<w:r w:rsidRPr="003331E2"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Escolha dos valores dos dados</w:t></w:r>'

Append-ConcluidoSuffix 11 $runs11
Append-ConcluidoSuffix 12 $runs12
Append-ConcluidoSuffix 13 $runs13
Append-ConcluidoSuffix 14 $runs14
Append-ConcluidoSuffix 15 $runs15

Write-Output "Done applying suffixes"
